$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.688358545303345
$ws.Range("B1").Value = 2.134612798690796
$ws.Range("C1").Value = 5.381175518035889
$ws.Range("D1").Value = 1.352210521697998
$ws.Range("E1").Value = 0.7536847591400146
